$d = $word.ActiveDocument

# Remove the placeholder text, leaving the bookmark intact.
$d.Content.Find.Execute("fegbaerugboaeitgboaeirtb", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# Append a single space after the (now-empty) paragraph content / bookmark.
$end = $d.Content.End
$r = $d.Range($end - 1, $end - 1)
$r.InsertAfter(" ")
